$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 1849.625
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 1849.625
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 1849.625
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -2199.625

# Row 38
$ws.Range("H38").Value = 1904.6
$ws.Range("J38").Value = 6575
$ws.Range("L38").Value = 19725
$ws.Range("N38").Value = -20469

# Row 40
$ws.Range("H40").Value = 3956.8
$ws.Range("I40").Value = 3740.5715
$ws.Range("J40").Value = 4461.3335
$ws.Range("K40").Value = 3740.5715
$ws.Range("L40").Value = 4461.3335
$ws.Range("M40").Value = -3565.5715
$ws.Range("N40").Value = -4811.3335

# Row 74
$ws.Range("H74").Value = 4202.1113
$ws.Range("J74").Value = 4496
$ws.Range("L74").Value = 4496
$ws.Range("N74").Value = -6368

# Row 77
$ws.Range("H77").Value = 4202.1113
$ws.Range("J77").Value = 4496
$ws.Range("L77").Value = 22480
$ws.Range("N77").Value = -31840

# Row 105
$ws.Range("H105").Value = 26999
$ws.Range("J105").Value = 26999
$ws.Range("L105").Value = 26999
$ws.Range("N105").Value = -33987

# Row 127
$ws.Range("H127").Value = 2882.2307
$ws.Range("J127").Value = 21700
$ws.Range("L127").Value = 65100
$ws.Range("N127").Value = -75020

# Row 132
$ws.Range("H132").Value = 3567.9143
$ws.Range("I132").Value = 3378.7354
$ws.Range("K132").Value = 10136.2062
$ws.Range("M132").Value = -7606.206200000001

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6956.5
$ws.Range("I32").Value = 7270.727
$ws.Range("K32").Value = 7270.727
$ws.Range("M32").Value = -6983.727

# Row 45
$ws.Range("H45").Value = 7353.636
$ws.Range("I45").Value = 6127.2856
$ws.Range("J45").Value = 9499.75
$ws.Range("K45").Value = 6127.2856
$ws.Range("L45").Value = 9499.75
$ws.Range("M45").Value = -5750.2856
$ws.Range("N45").Value = -10253.75

# Row 61
$ws.Range("H61").Value = 3973.5676
$ws.Range("I61").Value = 3361.36
$ws.Range("K61").Value = 3361.36
$ws.Range("M61").Value = -3149.36

# Row 122
$ws.Range("H122").Value = 271794.3
$ws.Range("I122").Value = 1823.7567
$ws.Range("J122").Value = 937721.7
$ws.Range("K122").Value = 5471.2701
$ws.Range("L122").Value = 2813165.1
$ws.Range("M122").Value = -3021.2701
$ws.Range("N122").Value = -2818065.1

# Row 132
$ws.Range("H132").Value = 3744.9333
$ws.Range("I132").Value = 2709.9473
$ws.Range("J132").Value = 5532.636
$ws.Range("K132").Value = 8129.841899999999
$ws.Range("L132").Value = 16597.908
$ws.Range("M132").Value = -5599.841899999999
$ws.Range("N132").Value = -21657.908

# Row 136
$ws.Range("H136").Value = 3973.5676
$ws.Range("I136").Value = 3361.36
$ws.Range("K136").Value = 10084.08
$ws.Range("M136").Value = -7534.08

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 4500.875
$ws.Range("I20").Value = 902.6667
$ws.Range("K20").Value = 902.6667
$ws.Range("M20").Value = -655.6667

# Row 64
$ws.Range("H64").Value = 7404.4287
$ws.Range("J64").Value = 3437.25
$ws.Range("L64").Value = 3437.25
$ws.Range("N64").Value = -3887.25

# Row 67
$ws.Range("H67").Value = 7404.4287
$ws.Range("J67").Value = 3437.25
$ws.Range("L67").Value = 3437.25
$ws.Range("N67").Value = -4997.25

# Row 86
$ws.Range("H86").Value = 3660.353
$ws.Range("I86").Value = 4355.8096
$ws.Range("J86").Value = 2536.923
$ws.Range("K86").Value = 4355.8096
$ws.Range("L86").Value = 2536.923
$ws.Range("M86").Value = -3232.8096
$ws.Range("N86").Value = -4782.923

# Row 89
$ws.Range("H89").Value = 3660.353
$ws.Range("I89").Value = 4355.8096
$ws.Range("J89").Value = 2536.923
$ws.Range("K89").Value = 21779.048
$ws.Range("L89").Value = 12684.615
$ws.Range("M89").Value = -16163.048
$ws.Range("N89").Value = -23916.615

# Row 99
$ws.Range("H99").Value = 27525.5
$ws.Range("I99").Value = 34528.184
$ws.Range("J99").Value = 1849
$ws.Range("K99").Value = 34528.184
$ws.Range("L99").Value = 1849
$ws.Range("M99").Value = -33030.184
$ws.Range("N99").Value = -4845

# Row 103
$ws.Range("H103").Value = 42850
$ws.Range("J103").Value = 42850
$ws.Range("L103").Value = 42850
$ws.Range("N103").Value = -45194

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5499.4287
$ws.Range("I31").Value = 4449
$ws.Range("K31").Value = 4449
$ws.Range("M31").Value = -4154

# Row 34
$ws.Range("H34").Value = 5499.4287
$ws.Range("I34").Value = 4449
$ws.Range("K34").Value = 4449
$ws.Range("M34").Value = -4247

# Row 52
$ws.Range("H52").Value = 46988.668
$ws.Range("J52").Value = 46988.668
$ws.Range("L52").Value = 46988.668
$ws.Range("N52").Value = -47576.668

# Row 58
$ws.Range("H58").Value = 2578.2285
$ws.Range("I58").Value = 1806.9166
$ws.Range("J58").Value = 4261.091
$ws.Range("K58").Value = 1806.9166
$ws.Range("L58").Value = 4261.091
$ws.Range("M58").Value = -1603.9166
$ws.Range("N58").Value = -4667.091

# Row 105
$ws.Range("H105").Value = 18616.5
$ws.Range("I105").Value = 21339.8
$ws.Range("K105").Value = 21339.8
$ws.Range("M105").Value = -19592.8

# Row 122
$ws.Range("H122").Value = 2031.7693
$ws.Range("I122").Value = 1899.8572
$ws.Range("J122").Value = 2185.6667
$ws.Range("K122").Value = 5699.571599999999
$ws.Range("L122").Value = 6557.000100000001
$ws.Range("M122").Value = -3249.571599999999
$ws.Range("N122").Value = -11457.0001

# Row 136
$ws.Range("H136").Value = 2578.2285
$ws.Range("I136").Value = 1806.9166
$ws.Range("J136").Value = 4261.091
$ws.Range("K136").Value = 5420.7498
$ws.Range("L136").Value = 12783.273
$ws.Range("M136").Value = -2870.7498
$ws.Range("N136").Value = -17883.273

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 599.6667
$ws.Range("I3").Value = 599.6667
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1799.0001
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1687.0001
$ws.Range("N3").ClearContents()

# Row 7
$ws.Range("H7").Value = 27.285715
$ws.Range("J7").Value = 6.3333335
$ws.Range("L7").Value = 19.0000005
$ws.Range("N7").Value = -243.0000005

# Row 34
$ws.Range("H34").Value = 1192761.1
$ws.Range("J34").Value = 1934
$ws.Range("L34").Value = 5802
$ws.Range("N34").Value = -5970

# Row 39
$ws.Range("H39").Value = 2138
$ws.Range("J39").Value = 4316.1665
$ws.Range("L39").Value = 12948.4995
$ws.Range("N39").Value = -13536.4995

# Row 55
$ws.Range("H55").Value = 5703.9287
$ws.Range("J55").Value = 6371.9165
$ws.Range("L55").Value = 19115.7495
$ws.Range("N55").Value = -19469.7495

# Row 63
$ws.Range("H63").Value = 2993.5

# Row 66
$ws.Range("H66").Value = 2993.5

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 7606.136
$ws.Range("I102").Value = 9333.4375
$ws.Range("K102").Value = 9333.4375
$ws.Range("M102").Value = -7711.4375

# Row 122
$ws.Range("H122").Value = 14941.5
$ws.Range("I122").Value = 14182.0625
$ws.Range("J122").Value = 16966.666
$ws.Range("K122").Value = 42546.1875
$ws.Range("L122").Value = 50899.99800000001
$ws.Range("M122").Value = -40096.1875
$ws.Range("N122").Value = -55799.99800000001

# Row 132
$ws.Range("H132").Value = 3007.923
$ws.Range("I132").Value = 3007.923
$ws.Range("K132").Value = 9023.769
$ws.Range("M132").Value = -6493.769

# Row 136
$ws.Range("H136").Value = 103577.89
$ws.Range("J136").Value = 103577.89
$ws.Range("L136").Value = 310733.67
$ws.Range("N136").Value = -315833.67

$ws = $wb.Worksheets.Item("LTW")
# Row 34
$ws.Range("H34").Value = 12755
$ws.Range("I34").Value = 13673.667
$ws.Range("J34").Value = 9999
$ws.Range("K34").Value = 13673.667
$ws.Range("L34").Value = 9999
$ws.Range("M34").Value = -13501.667
$ws.Range("N34").Value = -10343

# Row 46
$ws.Range("H46").Value = 4003.7144
$ws.Range("I46").Value = 1498.4286
$ws.Range("J46").Value = 5256.357
$ws.Range("K46").Value = 1498.4286
$ws.Range("L46").Value = 5256.357
$ws.Range("M46").Value = -1310.4286
$ws.Range("N46").Value = -5632.357

# Row 55
$ws.Range("H55").Value = 946
$ws.Range("I55").Value = 946
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 946
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -773
$ws.Range("N55").ClearContents()

# Row 93
$ws.Range("H93").Value = 5714.9033
$ws.Range("I93").Value = 6270.1924
$ws.Range("K93").Value = 6270.1924
$ws.Range("M93").Value = -5022.1924

# Row 97
$ws.Range("H97").Value = 15000
$ws.Range("J97").Value = 15000
$ws.Range("L97").Value = 15000
$ws.Range("N97").Value = -16982

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 22398.2
$ws.Range("J107").Value = 61990
$ws.Range("L107").Value = 185970
$ws.Range("N107").Value = -189810

# Row 122
$ws.Range("H122").Value = 6600.467
$ws.Range("I122").Value = 3438.4375
$ws.Range("K122").Value = 10315.3125
$ws.Range("M122").Value = -7865.3125

# Row 137
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
